$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.386.24"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.331.06"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.95"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.66"
$ws.Range("E6").Value = "  +4.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  +1.14%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.614"
$ws.Range("E9").Value = "  +2.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.64"
$ws.Range("E10").Value = "  +3.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0921"
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.55"
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.01"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.44"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.680.71"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.320.08"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.276.67"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.15"
$ws.Range("E21").Value = "  -10.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.06"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.50"
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.58"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.27"
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.61"
$ws.Range("E27").Value = "  +10.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.12"
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.15"
$ws.Range("E30").Value = "  +4.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.62"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.18"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0885"
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.77"
$ws.Range("E34").Value = "  +7.56%  "
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("E36").Value = "  +3.22%  "
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0361"
$ws.Range("E38").Value = "  +3.11%  "
$ws.Range("E39").Value = "  +5.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.76"
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.66"
$ws.Range("E41").Value = "  +5.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.95"
$ws.Range("E42").Value = "  +12.00%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.42"
$ws.Range("E43").Value = "  +9.36%  "
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.59"
$ws.Range("E44").Value = "  +3.21%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.235"
$ws.Range("E45").Value = "  +2.86%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "113.86"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.667.37"
$ws.Range("E48").Value = "  -3.88%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "77.58"
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.218"
$ws.Range("E50").Value = "  +17.09%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.96"
$ws.Range("E51").Value = "  +5.24%  "
